{"js": "// Remove the stray \"AND/OR\" text left over at the end of the document,\n// leaving the (now empty) paragraph itself in place.\nconst body = context.document.body;\nconst results = body.search(\"AND/OR\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (const result of results.items) {\n  result.delete();\n}\nawait context.sync();\n", "ps1": "# Remove the stray \"AND/OR\" text left over at the end of the document,\n# leaving the (now empty) paragraph itself in place.\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    $text = $r.Text.TrimEnd([char]13, [char]7)\n    if ($text -eq \"AND/OR\") {\n        $r.MoveEnd(1, -1) | Out-Null   # wdCharacter = 1; exclude the trailing paragraph mark\n        $r.Delete()\n    }\n}\n"}
